# Applies the "Add files via upload" revision to lit/Book1.xlsx:
#  - adds two new reference rows (7 & 8: Malta 2017 / 2018 SBA Fact Sheets)
#  - retitles a few existing rows (SBA Fact Sheet / SBA Fact Sheet and Scoreboard /
#    Key-success-factors... with the leading space trimmed)
#  - drops the now-unused "Year" (B) values/column usage on most rows and merges
#    the old two-column "X" marker (E+F) into a single E column
#  - normalizes the "X" marker cells' alignment (drop vertical=top, keep horizontal=center)
#  - gives the literature-review title cell (C6) a wrap-text / vertical-centered style

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------------------
$ws.Range("B2").Clear()
$ws.Range("E2").Clear()
$ws.Range("F2").Clear()

$ws.Range("G2").ClearFormats()
$ws.Range("H2").ClearFormats()
$ws.Range("I2").ClearFormats()
$ws.Range("G2").HorizontalAlignment = $xlCenter
$ws.Range("H2").HorizontalAlignment = $xlCenter
$ws.Range("I2").Value = "X"
$ws.Range("I2").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------------------
$ws.Range("B3").Clear()
$ws.Range("E3").Clear()
$ws.Range("F3").Clear()

$ws.Range("G3").ClearFormats()
$ws.Range("H3").ClearFormats()
$ws.Range("I3").ClearFormats()
$ws.Range("G3").HorizontalAlignment = $xlCenter
$ws.Range("H3").HorizontalAlignment = $xlCenter
$ws.Range("I3").Value = "X"
$ws.Range("I3").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = "SBA Fact Sheet"

$ws.Range("F4").Clear()

$ws.Range("E4").ClearFormats()
$ws.Range("E4").Value = "X"
$ws.Range("E4").HorizontalAlignment = $xlCenter
$ws.Range("E4").VerticalAlignment = $xlCenter

$ws.Range("G4").ClearFormats()
$ws.Range("G4").HorizontalAlignment = $xlCenter

$ws.Range("H4").ClearFormats()
$ws.Range("H4").ClearContents()
$ws.Range("H4").HorizontalAlignment = $xlCenter

$ws.Range("I4").ClearFormats()
$ws.Range("I4").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "SBA Fact Sheet and Scoreboard"

$ws.Range("F5").Clear()

$ws.Range("E5").ClearFormats()
$ws.Range("E5").Value = "X"
$ws.Range("E5").HorizontalAlignment = $xlCenter

$ws.Range("G5").ClearFormats()
$ws.Range("G5").HorizontalAlignment = $xlCenter

$ws.Range("H5").ClearFormats()
$ws.Range("H5").ClearContents()
$ws.Range("H5").HorizontalAlignment = $xlCenter

$ws.Range("I5").ClearFormats()
$ws.Range("I5").ClearContents()
$ws.Range("I5").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 6
# ---------------------------------------------------------------------------
$ws.Range("B6").Clear()
$ws.Range("E6").Clear()
$ws.Range("F6").Clear()

$ws.Range("C6").Value = "Key-success-factors-for-business-incubators-in-europe-an-empirical"
$ws.Range("C6").ClearFormats()
$ws.Range("C6").WrapText = $true
$ws.Range("C6").VerticalAlignment = $xlCenter

$ws.Range("G6").ClearFormats()
$ws.Range("H6").ClearFormats()
$ws.Range("I6").ClearFormats()
$ws.Range("G6").HorizontalAlignment = $xlCenter
$ws.Range("H6").HorizontalAlignment = $xlCenter
$ws.Range("I6").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 7 (new)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("C7").Value = "Malta - 2017 SBA Fact Sheet"
$ws.Range("E7").ClearFormats()
$ws.Range("E7").Value = "X"
$ws.Range("E7").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Row 8 (new)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 7
$ws.Range("C8").Value = "Malta - 2018 Fact Sheet"
$ws.Range("E8").ClearFormats()
$ws.Range("E8").Value = "X"
$ws.Range("E8").HorizontalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Selection / view tidy-up
# ---------------------------------------------------------------------------
$ws.Range("C14").Select()
